$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.186842249490837
$ws.Cells.Item(2, 4).Value = 0.007892881942677121
$ws.Cells.Item(2, 5).Value = 0.7787286082367899
$ws.Cells.Item(2, 6).Value = 0.4994444187775429
$ws.Cells.Item(2, 7).Value = 0.002369686683785481
$ws.Cells.Item(2, 9).Value = 0.9151550329007279
$ws.Cells.Item(2, 12).Value = 0.3317526960658199
$ws.Cells.Item(2, 13).Value = 0.3235158654437384
$ws.Cells.Item(2, 15).Value = 1.576574840463593
$ws.Cells.Item(3, 2).Value = 1.087854672962521
$ws.Cells.Item(3, 4).Value = 0.007091353438440251
$ws.Cells.Item(3, 5).Value = 0.7311082958859743
$ws.Cells.Item(3, 6).Value = 0.4862820119347688
$ws.Cells.Item(3, 7).Value = 0.002373058442449144
$ws.Cells.Item(3, 9).Value = 0.9151100917641557
$ws.Cells.Item(3, 12).Value = 0.298675550832229
$ws.Cells.Item(3, 13).Value = 0.2944447077293546
$ws.Cells.Item(3, 15).Value = 1.546693487707159
$ws.Cells.Item(4, 2).Value = 1.026938595950412
$ws.Cells.Item(4, 4).Value = 0.006597748167170181
$ws.Cells.Item(4, 5).Value = 0.7019212379060491
$ws.Cells.Item(4, 6).Value = 0.4786915913306444
$ws.Cells.Item(4, 7).Value = 0.002375239564484484
$ws.Cells.Item(4, 9).Value = 0.9158644192635634
$ws.Cells.Item(4, 12).Value = 0.2783327030838478
$ws.Cells.Item(4, 13).Value = 0.2765600652531717
$ws.Cells.Item(4, 15).Value = 1.529995555072077
$ws.Cells.Item(5, 2).Value = 1.00208187030745
$ws.Cells.Item(5, 4).Value = 0.006396250649228108
$ws.Cells.Item(5, 5).Value = 0.6900420218302941
$ws.Cells.Item(5, 6).Value = 0.475721362181261
$ws.Cells.Item(5, 7).Value = 0.002376156349263332
$ws.Cells.Item(5, 9).Value = 0.9163687134598391
$ws.Cells.Item(5, 12).Value = 0.2700349587544792
$ws.Cells.Item(5, 13).Value = 0.2692636462573219
$ws.Cells.Item(5, 15).Value = 1.523603524467205
$ws.Cells.Item(6, 2).Value = 0.9979524886025501
$ws.Cells.Item(6, 4).Value = 0.006362771592144156
$ws.Cells.Item(6, 5).Value = 0.6880704270949565
$ws.Cells.Item(6, 6).Value = 0.475235564179556
$ws.Cells.Item(6, 7).Value = 0.002376310272053409
$ws.Cells.Item(6, 9).Value = 0.9164643518848692
$ws.Cells.Item(6, 12).Value = 0.2686566647285815
$ws.Cells.Item(6, 13).Value = 0.2680515954187044
$ws.Cells.Item(6, 15).Value = 1.522566982132702
$ws.Cells.Item(7, 2).Value = 1.026603500655654
$ws.Cells.Item(7, 4).Value = 0.006595032091020414
$ws.Cells.Item(7, 5).Value = 0.7017609687081432
$ws.Cells.Item(7, 6).Value = 0.478651036878091
$ws.Cells.Item(7, 7).Value = 0.002375251815172628
$ws.Cells.Item(7, 9).Value = 0.9158704227331143
$ws.Cells.Item(7, 12).Value = 0.278220827972504
$ws.Cells.Item(7, 13).Value = 0.2764616960913386
$ws.Cells.Item(7, 15).Value = 1.529907682624781
$ws.Cells.Item(8, 2).Value = 1.152740750405144
$ws.Cells.Item(8, 4).Value = 0.007616829544364379
$ws.Cells.Item(8, 5).Value = 0.7622993121651405
$ws.Cells.Item(8, 6).Value = 0.494803660006383
$ws.Cells.Item(8, 7).Value = 0.002370826311746992
$ws.Cells.Item(8, 9).Value = 0.9149774079401851
$ws.Cells.Item(8, 12).Value = 0.3203549152076164
$ws.Cells.Item(8, 13).Value = 0.313499635440138
$ws.Cells.Item(8, 15).Value = 1.565928077472961
$ws.Cells.Item(9, 2).Value = 1.398948943771927
$ws.Cells.Item(9, 4).Value = 0.009608243528099081
$ws.Cells.Item(9, 5).Value = 0.8813605395253745
$ws.Cells.Item(9, 6).Value = 0.5304080247617691
$ws.Cells.Item(9, 7).Value = 0.002363023507560454
$ws.Cells.Item(9, 9).Value = 0.9194208113321238
$ws.Cells.Item(9, 12).Value = 0.4026967341413297
$ws.Cells.Item(9, 13).Value = 0.3858373529630086
$ws.Cells.Item(9, 15).Value = 1.649756960998076
$ws.Cells.Item(10, 2).Value = 1.579081531568249
$ws.Cells.Item(10, 4).Value = 0.0110630121429125
$ws.Cells.Item(10, 5).Value = 0.9689650480906806
$ws.Cells.Item(10, 6).Value = 0.5590073065671533
$ws.Cells.Item(10, 7).Value = 0.002357819051400917
$ws.Cells.Item(10, 9).Value = 0.926452806400313
$ws.Cells.Item(10, 12).Value = 0.4630023270622985
$ws.Cells.Item(10, 13).Value = 0.4387874069518887
$ws.Cells.Item(10, 15).Value = 1.719544211235814
$ws.Cells.Item(11, 2).Value = 1.660853610973902
$ws.Cells.Item(11, 4).Value = 0.01172285861587596
$ws.Cells.Item(11, 5).Value = 1.00883061113467
$ws.Cells.Item(11, 6).Value = 0.5725577905984096
$ws.Cells.Item(11, 7).Value = 0.002355564940847665
$ws.Cells.Item(11, 9).Value = 0.930468522341485
$ws.Cells.Item(11, 12).Value = 0.4903920343973596
$ws.Cells.Item(11, 13).Value = 0.4628298381722402
$ws.Cells.Item(11, 15).Value = 1.753106220837111
$ws.Cells.Item(12, 2).Value = 1.691792708571882
$ws.Cells.Item(12, 4).Value = 0.01197243104655143
$ws.Cells.Item(12, 5).Value = 1.023927126535355
$ws.Cells.Item(12, 6).Value = 0.5777674615905397
$ws.Cells.Item(12, 7).Value = 0.002354727588586868
$ws.Cells.Item(12, 9).Value = 0.932106482579286
$ws.Cells.Item(12, 12).Value = 0.5007571229129155
$ws.Cells.Item(12, 13).Value = 0.471927256485742
$ws.Cells.Item(12, 15).Value = 1.766078882290998
$ws.Cells.Item(13, 2).Value = 1.685130611471891
$ws.Cells.Item(13, 4).Value = 0.01191869466948248
$ws.Cells.Item(13, 5).Value = 1.020675835692614
$ws.Cells.Item(13, 6).Value = 0.5766419683661326
$ws.Cells.Item(13, 7).Value = 0.002354907206919616
$ws.Cells.Item(13, 9).Value = 0.9317485042073059
$ws.Cells.Item(13, 12).Value = 0.4985251263740622
$ws.Cells.Item(13, 13).Value = 0.469968280464073
$ws.Cells.Item(13, 15).Value = 1.763273232072464
$ws.Cells.Item(14, 2).Value = 1.663399523619034
$ws.Cells.Item(14, 4).Value = 0.01174339715146289
$ws.Cells.Item(14, 5).Value = 1.010072614352993
$ws.Cells.Item(14, 6).Value = 0.5729848183484876
$ws.Cells.Item(14, 7).Value = 0.002355495726441729
$ws.Cells.Item(14, 9).Value = 0.9306009284135115
$ws.Cells.Item(14, 12).Value = 0.4912449157334038
$ws.Cells.Item(14, 13).Value = 0.4635784298601209
$ws.Cells.Item(14, 15).Value = 1.754168195711742
$ws.Cells.Item(15, 2).Value = 1.650085143821457
$ws.Cells.Item(15, 4).Value = 0.01163598311343605
$ws.Cells.Item(15, 5).Value = 1.003577823054655
$ws.Cells.Item(15, 6).Value = 0.5707549377837893
$ws.Cells.Item(15, 7).Value = 0.002355858323119324
$ws.Cells.Item(15, 9).Value = 0.929913275899068
$ws.Cells.Item(15, 12).Value = 0.4867846754954712
$ws.Cells.Item(15, 13).Value = 0.4596635453025755
$ws.Cells.Item(15, 15).Value = 1.748625479007046
$ws.Cells.Item(16, 2).Value = 1.573733952320026
$ws.Cells.Item(16, 4).Value = 0.01101984906919284
$ws.Cells.Item(16, 5).Value = 0.966359872657037
$ws.Cells.Item(16, 6).Value = 0.5581326846368029
$ws.Cells.Item(16, 7).Value = 0.002357968638380523
$ws.Cells.Item(16, 9).Value = 0.9262067922850861
$ws.Cells.Item(16, 12).Value = 0.4612114224135837
$ws.Cells.Item(16, 13).Value = 0.437215234180627
$ws.Cells.Item(16, 15).Value = 1.717387591425421
$ws.Cells.Item(17, 2).Value = 1.526849986140178
$ws.Cells.Item(17, 4).Value = 0.01064136167510554
$ws.Cells.Item(17, 5).Value = 0.9435301953131301
$ws.Cells.Item(17, 6).Value = 0.5505282601459811
$ws.Cells.Item(17, 7).Value = 0.002359292241925086
$ws.Cells.Item(17, 9).Value = 0.9241420727551102
$ws.Cells.Item(17, 12).Value = 0.4455115123199391
$ws.Cells.Item(17, 13).Value = 0.4234320974733379
$ws.Cells.Item(17, 15).Value = 1.698690845098611
$ws.Cells.Item(18, 2).Value = 1.499867554621005
$ws.Cells.Item(18, 4).Value = 0.01042348470822674
$ws.Cells.Item(18, 5).Value = 0.9304005710824157
$ws.Cells.Item(18, 6).Value = 0.5462052333712535
$ws.Cells.Item(18, 7).Value = 0.002360064224210575
$ws.Cells.Item(18, 9).Value = 0.9230313929966201
$ws.Cells.Item(18, 12).Value = 0.4364772656597324
$ws.Cells.Item(18, 13).Value = 0.4155002121737823
$ws.Cells.Item(18, 15).Value = 1.688107705760473
$ws.Cells.Item(19, 2).Value = 1.490729065829726
$ws.Cells.Item(19, 4).Value = 0.01034968479803666
$ws.Cells.Item(19, 5).Value = 0.9259553987711797
$ws.Cells.Item(19, 6).Value = 0.5447502416715935
$ws.Cells.Item(19, 7).Value = 0.002360327441042487
$ws.Cells.Item(19, 9).Value = 0.9226685474537177
$ws.Cells.Item(19, 12).Value = 0.4334177415410068
$ws.Cells.Item(19, 13).Value = 0.4128139074114543
$ws.Cells.Item(19, 15).Value = 1.684553687302497
$ws.Cells.Item(20, 2).Value = 1.531842533497297
$ws.Cells.Item(20, 4).Value = 0.01068167117739449
$ws.Cells.Item(20, 5).Value = 0.9459603206521763
$ws.Cells.Item(20, 6).Value = 0.551332497812993
$ws.Cells.Item(20, 7).Value = 0.002359150237110501
$ws.Cells.Item(20, 9).Value = 0.9243539088346537
$ws.Cells.Item(20, 12).Value = 0.4471832200065364
$ws.Cells.Item(20, 13).Value = 0.4248997730470307
$ws.Cells.Item(20, 15).Value = 1.700663457872736
$ws.Cells.Item(21, 2).Value = 1.669783194987644
$ws.Cells.Item(21, 4).Value = 0.01179489449423698
$ws.Cells.Item(21, 5).Value = 1.013187043005516
$ws.Cells.Item(21, 6).Value = 0.5740568789649245
$ws.Cells.Item(21, 7).Value = 0.002355322424103035
$ws.Cells.Item(21, 9).Value = 0.9309348170916678
$ws.Cells.Item(21, 12).Value = 0.4933834787748594
$ws.Cells.Item(21, 13).Value = 0.4654554757389775
$ws.Cells.Item(21, 15).Value = 1.756835396362732
$ws.Cells.Item(22, 2).Value = 1.759781869801259
$ws.Cells.Item(22, 4).Value = 0.01252070995609245
$ws.Cells.Item(22, 5).Value = 1.057124518158517
$ws.Cells.Item(22, 6).Value = 0.5893658417036818
$ws.Cells.Item(22, 7).Value = 0.002352915294691991
$ws.Cells.Item(22, 9).Value = 0.9359194455762747
$ws.Cells.Item(22, 12).Value = 0.5235381477254748
$ws.Cells.Item(22, 13).Value = 0.4919203670962133
$ws.Cells.Item(22, 15).Value = 1.795083743336988
$ws.Cells.Item(23, 2).Value = 1.71176245358771
$ws.Cells.Item(23, 4).Value = 0.01213349398175012
$ws.Cells.Item(23, 5).Value = 1.033674732356374
$ws.Cells.Item(23, 6).Value = 0.5811530994798062
$ws.Cells.Item(23, 7).Value = 0.002354191398058189
$ws.Cells.Item(23, 9).Value = 0.9331965537751614
$ws.Cells.Item(23, 12).Value = 0.5074478541866085
$ws.Cells.Item(23, 13).Value = 0.4777994300057742
$ws.Cells.Item(23, 15).Value = 1.774528477667786
$ws.Cells.Item(24, 2).Value = 1.529585490595821
$ws.Cells.Item(24, 4).Value = 0.01066344812034004
$ws.Cells.Item(24, 5).Value = 0.944861674997739
$ws.Cells.Item(24, 6).Value = 0.5509687498237668
$ws.Cells.Item(24, 7).Value = 0.002359214403132001
$ws.Cells.Item(24, 9).Value = 0.9242578998885236
$ws.Cells.Item(24, 12).Value = 0.4464274663374965
$ws.Cells.Item(24, 13).Value = 0.4242362611055199
$ws.Cells.Item(24, 15).Value = 1.699771123141829
$ws.Cells.Item(25, 2).Value = 1.332471915431654
$ws.Cells.Item(25, 4).Value = 0.009070917597099992
$ws.Cells.Item(25, 5).Value = 0.849122703010778
$ws.Cells.Item(25, 6).Value = 0.5203507565200454
$ws.Cells.Item(25, 7).Value = 0.002365041207140476
$ws.Cells.Item(25, 9).Value = 0.9175565652032418
$ws.Cells.Item(25, 12).Value = 0.3804533207567431
$ws.Cells.Item(25, 13).Value = 0.3663013504448429
$ws.Cells.Item(25, 15).Value = 1.763273232072464
